$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- Row 4 content (Reseau Ongules sauvages) ---
$ws.Range("B4").Value = "Réseau Ongulés sauvages"
$ws.Range("D4").Value = "Le réseau a pour but de récolter les informations utiles pour suivre les 14 espèces d'ongulés sauvages présentes en France hexagonale"
$ws.Range("E4").Value = "L’ensemble des données récoltées permet de mesurer de façon régulière des variables biologiques d’intérêt, telles que les aires de présence, les tendances d’évolution ou l’état de santé des populations d’ongulés sauvages en France métropolitaine. En parallèle, d’autres données connexes, telles les prélèvements cynégétiques ou les modalités de gestion des populations, sont régulièrement enregistrées."
$ws.Range("G4").Value = $ws.Range("G5").Value2
$ws.Range("J4").Value = $ws.Range("J5").Value2
$ws.Range("L4").Value = "Animation nationale: XXXX`nAnimation régionale: Samuel DEMBSKI`nCorrespondants départementaux:`nPPC:`n77:`n78-95:`n91`nCourriel du réseau: reseau.ongules-sauvages@ofb.gouv.fr"
$ws.Range("M4").Value = "Fédérations de chasse"
$ws.Range("Z4").Value = $ws.Range("Z6").Value2
$ws.Range("AA4").Value = "texte:Dataviz: Présence des ongulés sauvages en France;lien:https://professionnels.ofb.fr/fr/doc-dataviz/dataviz-presence-ongules-sauvages-en-France"
$ws.Range("AB4").Value = "texte:Fiches de synthèse des suivis;lien:https://professionnels.ofb.fr/fr/node/869"
$ws.Range("AC4").Value = "texte: Cartes de répartition;lien:https://carmen.carmencarto.fr/38/Ongules_sauvages.map#"
$ws.Range("AD4").Value = "texte:Lettre d'information;lien:https://professionnels.ofb.fr/fr/node/1281"
$ws.Range("AE4").Value = "texte: Page du réseau sur le portail technique;lien:https://professionnels.ofb.fr/node/1431"

# --- Row height ---
$ws.Rows.Item(4).RowHeight = 285

# --- Selection / view state ---
$ws.Range("J4").Select()

Write-Host "done"
